# Update the dSF (column F) values per the repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = 4
$ws.Range("F10").Value = 1
$ws.Range("F13").Value = -9
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = 0
